# Apply the "Doc structure and overview sample" edit to conditions_v11.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- Build the new shared-string insertion order so the sharedStrings pool  ---
# --- ends up with the same sequence used by the target workbook.           ---
# 1) Fees
$ws.Range("K4").Value = "Fees"
$ws.Range("K4").Font.Bold = $true

# 2) Contract conditions (title, replaces old "Contract")
$ws.Range("A2").Value = "Contract conditions"

# 3) Non-benefiary calls / 4) Weekdays / 5) HUF (first occurrence on row 5)
$ws.Range("K5").Value = "Non-benefiary calls"
$ws.Range("L5").Value = "Weekdays"
$ws.Range("O5").Value = "HUF"

# 6) Benefiary calls / 7) Everyday (row 7)
$ws.Range("K7").Value = "Benefiary calls"
$ws.Range("L7").Value = "Everyday"

# 8) 00:00-23:59 (row7) / 9) 16:00-07:59 (row6) / 10) 08:00-15:59 (row5)
$ws.Range("M7").Value = "00:00-23:59"
$ws.Range("M6").Value = "16:00-07:59"
$ws.Range("M5").Value = "08:00-15:59"

# Remaining duplicate cells (reuse already interned strings)
$ws.Range("K6").Value = "Non-benefiary calls"
$ws.Range("L6").Value = "Weekdays"
$ws.Range("O6").Value = "HUF"
$ws.Range("O7").Value = "HUF"

# Numeric fee amounts
$ws.Range("N5").Value = 30
$ws.Range("N6").Value = 10
$ws.Range("N7").Value = 0

# Re-assert existing labels so the shared-string pool keeps them (text identical,
# only their pool index shifts because "Contract" is no longer referenced).
# NOTE: single-quoted strings are used for the template placeholders below so
# PowerShell does not try to interpolate the "${...}" / "$var" sequences.
$ws.Range("A4").Value = "Contractor"
$ws.Range("B5").Value = "Name"
$ws.Range("C5").Value = '${ctx[''contract''].contractor.name}'
$ws.Range("B6").Value = "Birth date"
$ws.Range("C6").Value = '${ctx[''contract''].contractor.birthDate}'
$ws.Range("A8").Value = "Beneficiaries"
$ws.Range("B9").Value = '${beneficiary.name}'
$ws.Range("C9").Value = '$beneficiary.phoneNumber}'

# New column widths for the fee table (closest achievable to the best-fit
# widths Excel computed: K ~16.9 chars, L:M ~10.8 chars)
$ws.Columns.Item(11).ColumnWidth = 16.0
$ws.Columns.Item(12).ColumnWidth = 10.0
$ws.Columns.Item(13).ColumnWidth = 10.0

# Update comment text (author prefix + body)
$comment = $ws.Range("B9").Comment
$commentText = "Szerző:`njx:each(items=`"ctx['contract'].beneficiaries`", groupBy=`"beneficiary.name`", var=`"beneficiary`", lastCell=`"C9`")"
$comment.Text($commentText) | Out-Null

# Selection moves to K7 after the edits
$ws.Range("K7").Select() | Out-Null

# Page setup additions on the other (empty) sheets
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1
$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1
